$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '96.778.88'
$ws.Range("E2").Value = '  -1.21%  '

# Row 3
$ws.Range("D3").Value = '3.668.38'
$ws.Range("E3").Value = '  +1.87%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.16'
$ws.Range("E5").Value = '  -1.04%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.89'
$ws.Range("E6").Value = '  +11.31%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '658.51'
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.426'
$ws.Range("E8").Value = '  +1.51%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.09'
$ws.Range("E9").Value = '  +3.12%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("E10").Value = '  +0.03%  '

# Row 11
$ws.Range("D11").Value = '3.665.31'
$ws.Range("E11").Value = '  +1.90%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.61'
$ws.Range("E12").Value = '  +2.81%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.206'
$ws.Range("E13").Value = '  +0.60%  '

# Row 14
$ws.Range("E14").Value = '  +4.68%  '

# Row 15
$ws.Range("D15").Value = '4.351.61'
$ws.Range("E15").Value = '  +1.93%  '

# Row 16
$ws.Range("E16").Value = '  +4.55%  '

# Row 17
$ws.Range("D17").Value = '96.539.10'
$ws.Range("E17").Value = '  -1.53%  '

# Row 18
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.89'
$ws.Range("E18").Value = '  +2.88%  '

# Row 19
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.675.02'
$ws.Range("E19").Value = '  +2.31%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.77'
$ws.Range("E20").Value = '  +4.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.78'
$ws.Range("E21").Value = '  +0.40%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.529'
$ws.Range("E22").Value = '  +1.90%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '533.46'
$ws.Range("E23").Value = '  +3.43%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.52'
$ws.Range("E24").Value = '  +0.86%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.15'
$ws.Range("E25").Value = '  +4.53%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000206'
$ws.Range("E26").Value = '  +0.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.30'
$ws.Range("E27").Value = '  +1.21%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.49'
$ws.Range("E28").Value = '  +3.99%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.168'
$ws.Range("E29").Value = '  +5.68%  '

# Row 30
$ws.Range("E30").Value = '  +5.57%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.04'
$ws.Range("E31").Value = '  +1.07%  '

# Row 32
$ws.Range("E32").Value = '  +0.03%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.92'
$ws.Range("E33").Value = '  +17.47%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.186'
$ws.Range("E34").Value = '  +0.86%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '665.88'
$ws.Range("E35").Value = '  +7.00%  '

# Row 36
$ws.Range("E36").Value = '  -0.33%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '32.58'
$ws.Range("E37").Value = '  +2.71%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.597'
$ws.Range("E38").Value = '  +4.72%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.87'
$ws.Range("E39").Value = '  +0.13%  '

# Row 40
$ws.Range("E40").Value = '  +3.55%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.00'
$ws.Range("E41").Value = '  +1.19%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.56'
$ws.Range("E42").Value = '  +9.67%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.959'
$ws.Range("E43").Value = '  +3.90%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.72'
$ws.Range("E44").Value = '  +16.16%  '

# Row 45
$ws.Range("E45").Value = '  +0.05%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0461'
$ws.Range("E46").Value = '  +4.22%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.436'
$ws.Range("E47").Value = '  +11.77%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.78'
$ws.Range("E48").Value = '  +5.85%  '

# Row 49
$ws.Range("E49").Value = '  +0.87%  '

# Row 50
$ws.Range("E50").Value = '  +0.02%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.71'
$ws.Range("E51").Value = '  +2.45%  '
